$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.611.38"
$ws.Range("E2").Value = "  +1.91%  "

$ws.Range("D3").Value = "1.598.30"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0908"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").Value = "1.825.34"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "1.610.24"
$ws.Range("E14").Value = "  +2.37%  "

$ws.Range("D15").Value = "29.604.63"
$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.537"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("E20").Value = "  +2.41%  "

$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0477"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("D35").Value = "1.431.46"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.49%  "

$ws.Range("E41").Value = "  +2.79%  "

$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "54.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0492"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.801"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.988"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "1.737.30"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
